# Regenerate the s_vals data to filter save games.
# Updates columns B (TB), C (d2S), D (K), E (IP) and the computed G (sum)
# for each data row (rows 2-12). Column F (Win) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732;  G = 9.295990156953671 }
    3  = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 4.429675500412797 }
    4  = @{ B = 0.127881588408715;  C = 0.3127903958511391; D = 0.1575252929769615; E = 0.496779210170732;  G = 1.094976487407548 }
    5  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 }
    6  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732;  G = 9.295990156953671 }
    7  = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 3.645393585217082 }
    8  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732;  G = 9.295990156953671 }
    9  = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732;  G = 7.524616544037286 }
    10 = @{ B = 1.459612070389937;  C = 0.3127903958511391; D = 0.1575252929769615; E = 0.496779210170732;  G = 2.42670696938877  }
    11 = @{ B = 0.6753301551942219; C = 0.3127903958511391; D = 0.8054896365839992; E = 0.496779210170732;  G = 2.290389397800092 }
    12 = @{ B = 0.6753301551942219; C = 0.3127903958511391; D = 3.900430680208489;  E = 0.496779210170732;  G = 5.385330441424582 }
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
